# Data refresh for "13 days before election" state grand-totals pull.
#
# For each affected state row, this updates the refreshed source totals
# (total_requested_2016 in column B and total_returned_2016 in column C;
# row 3 also corrects total_requested_2020 in column D), then recomputes
# the derived diff/pctchg columns (F: diff_requested, G: pctchg_requested,
# H: diff_returned, I: pctchg_returned) from the refreshed totals, exactly
# as the rest of the sheet's static values were originally computed:
#   diff_requested    = total_requested_2020 - total_requested_2016
#   pctchg_requested  = ROUND(diff_requested / total_requested_2016 * 100, 2)
#   diff_returned     = total_returned_2020 - total_returned_2016
#   pctchg_returned   = ROUND(diff_returned / total_returned_2016 * 100, 2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B=36805; C=19890 }
    3 = @{ B=165505; C=165490; D=110889 }
    4 = @{ B=2325694; C=802526 }
    5 = @{ B=10923038; C=1821848 }
    6 = @{ B=3165002; C=575690 }
    8 = @{ B=19836; C=14042 }
    9 = @{ B=4132873; C=2443037 }
    10 = @{ B=1053153; C=949515 }
    11 = @{ B=503686; C=350934 }
    12 = @{ B=125610; C=90875 }
    13 = @{ B=773025; C=507809 }
    14 = @{ B=237880; C=115659 }
    15 = @{ B=202337 }
    16 = @{ B=172190; C=113416 }
    17 = @{ B=1035908; C=572522 }
    18 = @{ B=230102; C=230069 }
    19 = @{ B=320788; C=142525 }
    20 = @{ B=1074555; C=962344 }
    21 = @{ B=75103; C=40765 }
    22 = @{ B=181955; C=101943 }
    23 = @{ B=368885; C=202440 }
    24 = @{ B=237995; C=205550 }
    25 = @{ B=301053; C=254900 }
    26 = @{ B=1441757; C=800450 }
    28 = @{ B=328253; C=328183 }
    30 = @{ B=32363; C=11241 }
    32 = @{ B=838409; C=838382 }
    33 = @{ B=1416680; C=1383036 }
    34 = @{ B=998026; C=172500 }
    35 = @{ B=386189; C=244453 }
    36 = @{ B=59705; C=37585 }
    37 = @{ B=441148; C=379398 }
    38 = @{ B=31732; C=24549 }
    39 = @{ B=55746; C=38695 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }

    $ws.Range("F$row").Formula = "=D$row-B$row"
    $ws.Range("G$row").Formula = "=ROUND((D$row-B$row)/B$row*100,2)"
    $ws.Range("H$row").Formula = "=E$row-C$row"
    $ws.Range("I$row").Formula = "=ROUND((E$row-C$row)/C$row*100,2)"
}

# Flatten the recomputed diff/pctchg formulas down to plain static values,
# matching the plain-value layout used throughout the rest of the sheet.
foreach ($row in $data.Keys) {
    foreach ($col in @("F","G","H","I")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = $cell.Value2
    }
}
